# Daily attendance processing - 2025-11-23 13:50:39
# Normalizes the "Recorded By" (column G) entries so that multiple
# recorders listed for a session are ordered consistently (reversed
# from their previous recorded order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -ne $null -and [string]$val -ne "") {
        $parts = [string]$val -split ", "
        if ($parts.Count -gt 1) {
            $reversedParts = $parts[($parts.Count - 1)..0]
            $newVal = $reversedParts -join ", "
            if ($newVal -ne [string]$val) {
                $cell.Value = $newVal
            }
        }
    }
}
